# 55_Inst_Workspace_RETAIL.xlsx - "aggiornamento excel retail e test XRA_factoring"
#
# 1) Fix a typo in the Business Line workspace id (dash -> underscore):
#    WS_IT_RETAIL-BL_RETAIL_000001 -> WS_IT_RETAIL_BL_RETAIL_000001
# 2) Make column C the same (best-fit) width as column B on the
#    "r Workspace_BusinessLine" sheet.
# 3) Refresh the remembered cell selections on the "Workspace" and
#    "r Workspace_BusinessLine" sheets, leaving the latter as the active tab.

$wb = $excel.ActiveWorkbook

$wsWorkspace = $wb.Worksheets.Item("Workspace")
$wsBL        = $wb.Worksheets.Item("r Workspace_BusinessLine")

# --- 1) Correct the mistyped workspace/business-line id -------------------
$wsBL.Range("B3").Value = "WS_IT_RETAIL_BL_RETAIL_000001"
$wsBL.Range("C3").Value = "WS_IT_RETAIL_BL_RETAIL_000001"

# --- 2) Match column C's width to column B's on the BusinessLine sheet ----
$bestFitWidth = $wsBL.Columns.Item(2).ColumnWidth
$wsBL.Columns.Item(3).ColumnWidth = $bestFitWidth

# --- 3) Update the saved selections ---------------------------------------
# "Workspace" sheet: A3:B3 -> C3 (select it first so it doesn't stay the
# active tab once we activate the BusinessLine sheet below).
[void]$wsWorkspace.Range("C3").Select()

# "r Workspace_BusinessLine" sheet: C10 -> C11, and keep this sheet active
# (it was the active tab before the edit, and remains so after).
[void]$wsBL.Range("C11").Select()
